# Fix typo in Cheat sheet with regex
#
# The "regex (正規表現) の例" table on slide 1 (graphicFrame "Table 84")
# has a cell whose text accidentally starts with two single-quote
# characters instead of one: "''^(?!Species$).*'" should read
# "'^(?!Species$).*'".

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)

    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shp = $slide.Shapes.Item($shi)

        if ($shp.HasTable) {
            $tbl = $shp.Table

            for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
                for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
                    $cell = $tbl.Cell($r, $c)
                    $tr = $cell.Shape.TextFrame.TextRange
                    $txt = $tr.Text

                    if ($txt -eq "''^(?!Species`$).*'") {
                        $tr.Text = "'^(?!Species`$).*'"
                    }
                }
            }
        }
    }
}
